# daily auto push: 2026-01-08 09:39 UTC
# Insert a new data row at row 607 (a new "2026/01/08 17:00" sample),
# which shifts the existing rows 607-648 down to 608-649 and extends
# the sheet's used range to A1:D649.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 607..648 down to 608..649, opening up a blank row 607.
$ws.Rows.Item(607).Insert()

# Populate the newly opened row 607 with the new observation.
# Column A holds date-like text (e.g. "2026/01/08") that must stay a
# plain string, not get auto-converted into a real date serial value,
# so it is entered with a leading apostrophe (force-text) and then the
# resulting "Text" number-format override is cleared so the cell ends
# up unstyled, same as every other date cell in the column.
$ws.Range("A607").Value = "'2026/01/08"
$ws.Range("B607").Value = "木"
$ws.Range("C607").Value = 17
$ws.Range("D607").Value = 199
$ws.Range("A607:D607").ClearFormats()
